# modificacion de rutas relativas
# Adds three new game-report rows (6, 7, 8) to the sheet, mirroring the
# existing table structure (columns A:Y).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 ---------------------------------------------------------------
$ws.Range("A6").Value = "2024-08-01 15:23:05"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 6
$ws.Range("P6").Value = 3
$ws.Range("R6").Value = 5
$ws.Range("T6").Value = 20
$ws.Range("U6").Value = 0
$ws.Range("V6").Value = "C:\Users\jgcorrea\Desktop\Repositorio\Predictor_ruleta\Data\Electromecanica.xlsx"
$ws.Range("X6").Value = "No es Simulación"
$ws.Range("Y6").Value = 0

# --- Row 7 ---------------------------------------------------------------
$ws.Range("A7").Value = "2024-08-01 15:27:01"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 6
$ws.Range("P7").Value = 3
$ws.Range("R7").Value = 5
$ws.Range("T7").Value = 20
$ws.Range("U7").Value = 0.4
$ws.Range("V7").Value = "C:\Users\jgcorrea\Desktop\Repositorio\Predictor_ruleta\Data\Crupier.xlsx"
$ws.Range("X7").Value = "No es Simulación"
$ws.Range("Y7").Value = 5

# --- Row 8 ---------------------------------------------------------------
$ws.Range("A8").Value = "2024-08-01 15:33:43"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("N8").Value = 10
$ws.Range("O8").Value = 6
$ws.Range("P8").Value = 3
$ws.Range("R8").Value = 5
$ws.Range("T8").Value = 20
$ws.Range("U8").Value = 0
$ws.Range("V8").Value = "C:\Users\jgcorrea\Desktop\Repositorio\Predictor_ruleta\Data\Electromecanica.xlsx"
$ws.Range("X8").Value = "No es Simulación"
$ws.Range("Y8").Value = 0
